# feat: update MacOs installation
#
# Adds new MacOS app-installation entries to the "Mac installed" sheet
# (Package / Aktiv / Brew / Notes columns) and a "Wine UI (Steam)" source
# link at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mac installed")

# --- New data rows -------------------------------------------------------
# Column A cells are written first, in the same order the new unique
# strings were originally introduced, so the shared-string table lines up
# with the authored workbook; the remaining (already-known) values follow.

$ws.Range("A35").Value = "htop"
$ws.Range("A32").Value = "whisky"
$ws.Range("A33").Value = "Heroic Games Launcher"
$ws.Range("D33").Value = "https://heroicgameslauncher.com/downloads (Epic, GOG, Amazon)"
$ws.Range("D32").Value = "Wine UI (Steam)"
$ws.Range("A38").Value = "https://www.golem.de/news/tools-fuer-das-game-porting-toolkit-windows-games-auf-dem-mac-spielen-2401-180913-3.html"
$ws.Range("A34").Value = "Snap"
$ws.Range("D34").Value = "App Store"
$ws.Range("A36").Value = "pcloud"

$ws.Range("B32").Value = "✅"
$ws.Range("C32").Value = "✅"

$ws.Range("B33").Value = "✅"
$ws.Range("C33").Value = "?"

$ws.Range("B34").Value = "✅"
$ws.Range("C34").Value = "?"

$ws.Range("B35").Value = "✅"
$ws.Range("C35").Value = "✅"

$ws.Range("B36").Value = "✅"
$ws.Range("C36").Value = "❌"
$ws.Range("D36").Value = "Separate Download"

# --- Hyperlink for the golem.de source note at row 38 ---------------------
$ws.Hyperlinks.Add($ws.Range("A38"), "https://www.golem.de/news/tools-fuer-das-game-porting-toolkit-windows-games-auf-dem-mac-spielen-2401-180913-3.html")

# --- Restore selection near the newly added rows --------------------------
$ws.Activate()
$ws.Range("A36").Select()
